# Weekly fruit/vegetable price update: two new daily records (2022-08-03,
# Primera/Segunda, Brasil origin) were inserted into the "Mango" log right
# before the existing 2021-10-25 entries, pushing every subsequent row down
# by two positions (old row 392 -> new row 394, ..., old row 441 -> new row 443).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 392; everything from 392 downward shifts
# down by two rows (392->394 ... 441->443), and the sheet dimension grows
# from A1:T441 to A1:T443 automatically.
$ws.Rows.Item(392).Resize(2).Insert()

# New row 392: Femacal de La Calera, Coquimbo, 2022-08-03 (serial 44776),
# Mango, Primera, 228 units, 9000/9000/9000, Brasil, 2250 $/kg.
$ws.Range("A392").Value() = 3
$ws.Range("B392").Value() = "Femacal de La Calera"
$ws.Range("C392").Value() = "Coquimbo"
$ws.Range("D392").Value() = 44776
$ws.Range("E392").Value() = 5
$ws.Range("F392").Value() = "Fruta"
$ws.Range("G392").Value() = 100108
$ws.Range("H392").Value() = "Tropicales y subtropicales"
$ws.Range("I392").Value() = 100108002
$ws.Range("J392").Value() = "Mango"
$ws.Range("K392").Value() = "Sin especificar"
$ws.Range("L392").Value() = "Primera"
$ws.Range("M392").Value() = 228
$ws.Range("N392").Value() = 9000
$ws.Range("O392").Value() = 9000
$ws.Range("P392").Value() = 9000
$ws.Range("Q392").Value() = "`$/bandeja 4 kilos"
$ws.Range("R392").Value() = "Brasil"
$ws.Range("S392").Value() = 2250
$ws.Range("T392").Value() = 4

# New row 393: same date/lot, Segunda quality, same volumes/prices/origin.
$ws.Range("A393").Value() = 3
$ws.Range("B393").Value() = "Femacal de La Calera"
$ws.Range("C393").Value() = "Coquimbo"
$ws.Range("D393").Value() = 44776
$ws.Range("E393").Value() = 5
$ws.Range("F393").Value() = "Fruta"
$ws.Range("G393").Value() = 100108
$ws.Range("H393").Value() = "Tropicales y subtropicales"
$ws.Range("I393").Value() = 100108002
$ws.Range("J393").Value() = "Mango"
$ws.Range("K393").Value() = "Sin especificar"
$ws.Range("L393").Value() = "Segunda"
$ws.Range("M393").Value() = 228
$ws.Range("N393").Value() = 9000
$ws.Range("O393").Value() = 9000
$ws.Range("P393").Value() = 9000
$ws.Range("Q393").Value() = "`$/bandeja 4 kilos"
$ws.Range("R393").Value() = "Brasil"
$ws.Range("S393").Value() = 2250
$ws.Range("T393").Value() = 4
